$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42; this shifts existing rows 42-95 down to 43-96
# and automatically grows the sheet dimension from A1:R95 to A1:R96.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new record's data.
$ws.Range("A42").Value = 3
$ws.Range("B42").Value = "Femacal de La Calera"
$ws.Range("C42").Value = "Coquimbo"
$ws.Range("D42").Value = 44413
$ws.Range("E42").Value = 5
$ws.Range("F42").Value = 100112010
$ws.Range("G42").Value = "Achicoria"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 105
$ws.Range("K42").Value = 7000
$ws.Range("L42").Value = 7500
$ws.Range("M42").Value = 7262
$ws.Range("N42").Value = "$/caja 16 unidades"
$ws.Range("O42").Value = "Provincia de Quillota"
$ws.Range("P42").Value = 454
$ws.Range("Q42").Value = 16
$ws.Range("R42").Value = "Hortaliza"
